$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing columns B:D (RawActivations, PercActivations,
# totalActivation) one place to the right, into C:E, to make room for a new
# "segments" column in B. Cells are copied right-to-left (D->E, then C->D,
# then B->C) so that a source column is never overwritten before it has been
# read.
for ($row = 1; $row -le 20; $row++) {
    $ws.Cells.Item($row, 4).Copy($ws.Cells.Item($row, 5))
    $ws.Cells.Item($row, 3).Copy($ws.Cells.Item($row, 4))
    $ws.Cells.Item($row, 2).Copy($ws.Cells.Item($row, 3))
}

# New header for the freshly vacated column B; reuse the bold/bordered header
# style already used by the other headers (copy format from C1, which is the
# shifted "RawActivations" header).
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "segments"

# Move the segment-name labels that used to live in column A (rows 2-20) into
# the newly freed column B (clearing the header-style formatting that column
# A's data carried), and replace column A with a 0-based numeric index for
# each segment.
for ($row = 2; $row -le 20; $row++) {
    $name = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 3).Copy($ws.Cells.Item($row, 2))
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 1).Value = $row - 2
}
